$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 322
$ws1.Cells.Item(3, 6).Value = 1274
$ws1.Cells.Item(4, 6).Value = 368
$ws1.Cells.Item(5, 6).Value = 341
$ws1.Cells.Item(6, 6).Value = 3856
$ws1.Cells.Item(8, 6).Value = 756
$ws1.Cells.Item(9, 6).Value = 2248
$ws1.Cells.Item(11, 6).Value = 223
$ws1.Cells.Item(12, 6).Value = 744
$ws1.Cells.Item(13, 6).Value = 170
$ws1.Cells.Item(14, 6).Value = 168
$ws1.Cells.Item(15, 6).Value = 2169
$ws1.Cells.Item(17, 6).Value = 12
$ws1.Cells.Item(19, 6).Value = 339
$ws1.Cells.Item(21, 6).Value = 32

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(9, 6).Value = 96
$ws2.Cells.Item(11, 6).Value = 89
$ws2.Cells.Item(22, 6).Value = 57

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(4, 6).Value = 2089
$ws3.Cells.Item(5, 6).Value = 321

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 2089
$ws4.Cells.Item(5, 6).Value = 321
$ws4.Cells.Item(10, 6).Value = 322
$ws4.Cells.Item(11, 6).Value = 1274
$ws4.Cells.Item(12, 6).Value = 368
$ws4.Cells.Item(16, 6).Value = 341
$ws4.Cells.Item(17, 6).Value = 3856
$ws4.Cells.Item(20, 6).Value = 96
$ws4.Cells.Item(22, 6).Value = 89
$ws4.Cells.Item(23, 6).Value = 756
$ws4.Cells.Item(24, 6).Value = 2248
$ws4.Cells.Item(27, 6).Value = 223
$ws4.Cells.Item(28, 6).Value = 744
$ws4.Cells.Item(29, 6).Value = 170
$ws4.Cells.Item(30, 6).Value = 168
$ws4.Cells.Item(32, 6).Value = 2169
$ws4.Cells.Item(36, 6).Value = 12
$ws4.Cells.Item(38, 6).Value = 339
$ws4.Cells.Item(40, 6).Value = 32
$ws4.Cells.Item(47, 6).Value = 57
